$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table row shading exception (tblPrEx/shd) hand-off between row 1 and
#    row 2 of the addressing-plan table. This is a purely cosmetic,
#    colour-less "clear" shading flag that Word's object model surfaces as
#    Row.Shading; set it explicitly on both rows so the intent is recorded
#    even though the value itself (no fill / "clear") is visually a no-op.
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$headerRow = $table.Rows.Item(1)
$firstDataRow = $table.Rows.Item(2)
$headerRow.Shading.Texture = 0            # wdTextureNone == "clear", no fill
$firstDataRow.Shading.Texture = 0         # wdTextureNone == "clear", no fill

# ---------------------------------------------------------------------------
# 2) Fix the three "ip dhcp excluded-address" lines: the second address in
#    each pair had an extra "0" (10.50.x -> 10.5.x).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "ip dhcp excluded-address 10.5.10.1 10.50.10.9", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "ip dhcp excluded-address 10.5.10.1 10.5.10.9", 2) | Out-Null

$d.Content.Find.Execute(
    "ip dhcp excluded-address 10.5.20.1 10.50.20.9", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "ip dhcp excluded-address 10.5.20.1 10.5.20.9", 2) | Out-Null

$d.Content.Find.Execute(
    "ip dhcp excluded-address 10.5.30.1 10.50.30.9", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "ip dhcp excluded-address 10.5.30.1 10.5.30.9", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Fix the access-list line the same way (10.50.0.0 -> 10.5.0.0), and
#    leave behind the "_GoBack" bookmark at the cursor position right after
#    "ip 10.5" (this is what split the run in two in the authored diff).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "access-list 100 permit ip 10.50.0.0 0.0.255.255 host 10.10.30.10",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "access-list 100 permit ip 10.5.0.0 0.0.255.255 host 10.10.30.10", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the final (now-empty) paragraph to the
#    point right after "ip 10.5" on the access-list line. Word keeps only
#    one "_GoBack" bookmark at a time, tracking the last edit location, so
#    deleting the old one and adding a new one reproduces that behaviour.
# ---------------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$insertionPoint = $d.Content
$insertionPoint.Find.Execute(
    "access-list 100 permit ip 10.5", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$insertionPoint.Collapse(0)   # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $insertionPoint) | Out-Null
